$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the three data rows (2,3,4) of the "Artfynd" sheet:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# (same record set, just re-ordered - all other cell content is unchanged).
#
# We snapshot each row's values first (columns A:AY, the sheet's full used
# range) and then write the rotated snapshots back, so the three writes
# cannot clobber each other's source data.

$row2 = $ws.Range("A2:AY2").Value2
$row3 = $ws.Range("A3:AY3").Value2
$row4 = $ws.Range("A4:AY4").Value2

$ws.Range("A2:AY2").Value2 = $row4
$ws.Range("A3:AY3").Value2 = $row2
$ws.Range("A4:AY4").Value2 = $row3

# Columns Y (Startdatum), AA (Slutdatum) and I (Antal) store plain text that
# happens to look like a date / number ("1996-08-10", "2020-07-30", "1").
# A bulk .Value2 array write lets Excel's usual autoconvert reinterpret that
# text as a real date serial / number, which would change the stored type.
# Re-apply those specific cells as text explicitly, restoring the default
# "Normal" style afterwards so no stray per-cell number format lingers.
function Set-TextValue($range, $text) {
  $r = $ws.Range($range)
  $r.NumberFormat = "@"
  $r.Value2 = $text
  $r.Style = "Normal"
}

Set-TextValue "I2" "1"
Set-TextValue "I3" ""
Set-TextValue "I4" ""

Set-TextValue "Y2" "2020-07-30"
Set-TextValue "AA2" "2020-07-30"
Set-TextValue "Y3" "1996-08-10"
Set-TextValue "AA3" "1996-08-10"
Set-TextValue "Y4" "1996-08-10"
Set-TextValue "AA4" "1996-08-10"
